# Apply cell-level updates to the crypto price table (Coin/Link/Price/Volume)
# matching the scraped source-site refresh for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.261.91"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "2.472.01"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.93"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.60"
$ws.Range("E6").Value = "  -1.77%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("D9").Value = "2.472.48"
$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("E10").Value = "  -4.57%  "

$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  -3.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.83"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").Value = "2.924.93"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").Value = "68.886.21"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("E16").Value = "  -3.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.72"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").Value = "2.461.15"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.83"
$ws.Range("E19").Value = "  -3.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.37"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.10"
$ws.Range("E21").Value = "  -5.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.81"
$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  +1.34%  "

$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "67.20"
$ws.Range("E26").Value = "  -3.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.70"
$ws.Range("E27").Value = "  -2.80%  "

$ws.Range("D28").Value = "2.599.47"
$ws.Range("E28").Value = "  -0.84%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.25"
$ws.Range("E29").Value = "  -4.63%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").Value = "0.0₃0820"
$ws.Range("E31").Value = "  -6.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.22"
$ws.Range("E32").Value = "  -4.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "442.45"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.15"
$ws.Range("E35").Value = "  -4.20%  "

$ws.Range("E36").Value = "  -5.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.03"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("E40").Value = "  -3.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.95"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.303"
$ws.Range("E42").Value = "  -3.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.47"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.49"
$ws.Range("E45").Value = "  -5.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.10"
$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  -4.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.32"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("E49").Value = "  -2.14%  "

$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("E51").Value = "  -4.37%  "
